$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("updated_site_list")

# Copy formatting for the new row from existing rows that already carry the
# target styles (s=1 for A:E, s=2 for F, s=4 for G), then fill in the values.
$ws.Range("A9").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F9").Copy()
$ws.Range("F18").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A3").Copy()
$ws.Range("G18").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A18").Value = "CLP"
$ws.Range("B18").Value = "Cole Park"
$ws.Range("C18").Value = "Corpus Christi"
$ws.Range("D18").Value = "Corpus Christi Bay"
$ws.Range("E18").Value = "CCB"
$ws.Range("F18").Value = "27.776309, -97.391421"
$ws.Range("G18").Value = "At edge of park near marina by seawall"
